$p = $ppt.ActivePresentation

# 1. Change the table style on slide 16 (3rd shape, a graphicFrame/table)
$s = $p.Slides.Item(16)
$shape = $s.Shapes.Item(3)
$shape.Table.ApplyStyle("{46BF6B3A-1AA6-4572-B0A6-3627D734A84E}")

# 2. Re-point the presentation's theme colors ("Integral" -> "Office") on the
#    live theme part (the one the slide master / presentation actually use).
$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme
$officeRGB = @(
    0,          # dk1    000000
    16777215,   # lt1    FFFFFF
    6968388,    # dk2    44546A
    15132391,   # lt2    E7E6E6
    13998939,   # accent1 5B9BD5
    3243501,    # accent2 ED7D31
    10855845,   # accent3 A5A5A5
    49407,      # accent4 FFC000
    12874308,   # accent5 4472C4
    4697456,    # accent6 70AD47
    12673797,   # hlink   0563C1
    7491477     # folHlink 954F72
)
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeRGB[$i - 1]
}
